$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C80").Value = "Yes"
$ws.Range("C4:C80").Select()
$excel.ActiveWindow.ScrollRow = 78
$excel.ActiveWindow.ScrollColumn = 1
